$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$map = @(
    @{Row=1; Col=1; Old="18+10=28"; New="75-26=49"},
    @{Row=1; Col=2; Old="5+21=26"; New="75-13=62"},
    @{Row=1; Col=3; Old="11+53=64"; New="54+24=78"},
    @{Row=1; Col=4; Old="33+1=34"; New="89-73=16"},
    @{Row=1; Col=5; Old="10+64=74"; New="11-2=9"},
    @{Row=2; Col=1; Old="77-18=59"; New="58+11=69"},
    @{Row=2; Col=2; Old="45+4=49"; New="44-9=35"},
    @{Row=2; Col=3; Old="61-19=42"; New="29-20=9"},
    @{Row=2; Col=4; Old="18+30=48"; New="47+46=93"},
    @{Row=2; Col=5; Old="91+7=98"; New="30+7=37"},
    @{Row=3; Col=1; Old="5+81=86"; New="70-50=20"},
    @{Row=3; Col=2; Old="99-31=68"; New="32+35=67"},
    @{Row=3; Col=3; Old="5+36=41"; New="45-22=23"},
    @{Row=3; Col=4; Old="9+82=91"; New="28+3=31"},
    @{Row=3; Col=5; Old="84-59=25"; New="67+28=95"},
    @{Row=4; Col=1; Old="63-22=41"; New="0+48=48"},
    @{Row=4; Col=2; Old="1+59=60"; New="73-71=2"},
    @{Row=4; Col=3; Old="42+55=97"; New="58+4=62"},
    @{Row=4; Col=4; Old="58+7=65"; New="62-59=3"},
    @{Row=4; Col=5; Old="45+22=67"; New="4+89=93"},
    @{Row=5; Col=1; Old="5+33=38"; New="17+49=66"},
    @{Row=5; Col=2; Old="16+61=77"; New="81-15=66"},
    @{Row=5; Col=3; Old="99-58=41"; New="76-64=12"},
    @{Row=5; Col=4; Old="19+27=46"; New="46+4=50"},
    @{Row=5; Col=5; Old="41-39=2"; New="93-10=83"},
    @{Row=6; Col=1; Old="87-85=2"; New="8+12=20"},
    @{Row=6; Col=2; Old="19+21=40"; New="44+14=58"},
    @{Row=6; Col=3; Old="84-79=5"; New="65+8=73"},
    @{Row=6; Col=4; Old="44+30=74"; New="12+36=48"},
    @{Row=6; Col=5; Old="90+7=97"; New="44-40=4"},
    @{Row=7; Col=1; Old="89-21=68"; New="82-5=77"},
    @{Row=7; Col=2; Old="71-21=50"; New="20-6=14"},
    @{Row=7; Col=3; Old="99-60=39"; New="64-43=21"},
    @{Row=7; Col=4; Old="51+12=63"; New="69+13=82"},
    @{Row=7; Col=5; Old="95-50=45"; New="41-16=25"},
    @{Row=8; Col=1; Old="69+5=74"; New="78+11=89"},
    @{Row=8; Col=2; Old="4+22=26"; New="30+47=77"},
    @{Row=8; Col=3; Old="60-3=57"; New="87-13=74"},
    @{Row=8; Col=4; Old="55+37=92"; New="23+4=27"},
    @{Row=8; Col=5; Old="67-54=13"; New="34-13=21"},
    @{Row=9; Col=1; Old="67+18=85"; New="77-58=19"},
    @{Row=9; Col=2; Old="97-31=66"; New="90-67=23"},
    @{Row=9; Col=3; Old="21+46=67"; New="27-26=1"},
    @{Row=9; Col=4; Old="24+67=91"; New="0+11=11"},
    @{Row=9; Col=5; Old="4+10=14"; New="39+21=60"},
    @{Row=10; Col=1; Old="7+27=34"; New="50+29=79"},
    @{Row=10; Col=2; Old="45-4=41"; New="38-35=3"},
    @{Row=10; Col=3; Old="1+81=82"; New="39+26=65"},
    @{Row=10; Col=4; Old="35+1=36"; New="67-25=42"},
    @{Row=10; Col=5; Old="93-36=57"; New="21-7=14"},
    @{Row=11; Col=1; Old="80-75=5"; New="34+6=40"},
    @{Row=11; Col=2; Old="2+13=15"; New="16-2=14"},
    @{Row=11; Col=3; Old="93-83=10"; New="70+5=75"},
    @{Row=11; Col=4; Old="87+9=96"; New="64+26=90"},
    @{Row=11; Col=5; Old="73+16=89"; New="57+17=74"},
    @{Row=12; Col=1; Old="40-31=9"; New="25+23=48"},
    @{Row=12; Col=2; Old="66-52=14"; New="90+3=93"},
    @{Row=12; Col=3; Old="24-14=10"; New="8+83=91"},
    @{Row=12; Col=4; Old="74+1=75"; New="44-14=30"},
    @{Row=12; Col=5; Old="79-61=18"; New="62-42=20"},
    @{Row=13; Col=1; Old="55-43=12"; New="83-37=46"},
    @{Row=13; Col=2; Old="87-6=81"; New="44-10=34"},
    @{Row=13; Col=3; Old="31-31=0"; New="61-28=33"},
    @{Row=13; Col=4; Old="45+37=82"; New="28+68=96"},
    @{Row=13; Col=5; Old="44+10=54"; New="15+33=48"},
    @{Row=14; Col=1; Old="75+9=84"; New="17+58=75"},
    @{Row=14; Col=2; Old="35+17=52"; New="76-32=44"},
    @{Row=14; Col=3; Old="66+26=92"; New="22+63=85"},
    @{Row=14; Col=4; Old="9+47=56"; New="35+56=91"},
    @{Row=14; Col=5; Old="71-39=32"; New="84-14=70"},
    @{Row=15; Col=1; Old="69+24=93"; New="99-84=15"},
    @{Row=15; Col=2; Old="29-11=18"; New="31+62=93"},
    @{Row=15; Col=3; Old="93-18=75"; New="38-32=6"},
    @{Row=15; Col=4; Old="19+72=91"; New="20+71=91"},
    @{Row=15; Col=5; Old="20+65=85"; New="82-7=75"},
    @{Row=16; Col=1; Old="7+51=58"; New="37-18=19"},
    @{Row=16; Col=2; Old="34+17=51"; New="59+11=70"},
    @{Row=16; Col=3; Old="5-2=3"; New="99-38=61"},
    @{Row=16; Col=4; Old="15+45=60"; New="28+27=55"},
    @{Row=16; Col=5; Old="97-76=21"; New="49+32=81"},
    @{Row=17; Col=1; Old="65+18=83"; New="62+10=72"},
    @{Row=17; Col=2; Old="17+61=78"; New="6+42=48"},
    @{Row=17; Col=3; Old="8-5=3"; New="53-29=24"},
    @{Row=17; Col=4; Old="62+0=62"; New="37-17=20"},
    @{Row=17; Col=5; Old="58-36=22"; New="54-35=19"},
    @{Row=18; Col=1; Old="38+44=82"; New="94-59=35"},
    @{Row=18; Col=2; Old="70-11=59"; New="48-23=25"},
    @{Row=18; Col=3; Old="57-17=40"; New="66+31=97"},
    @{Row=18; Col=4; Old="67-30=37"; New="47+40=87"},
    @{Row=18; Col=5; Old="41+48=89"; New="68+15=83"},
    @{Row=19; Col=1; Old="96-66=30"; New="27+9=36"},
    @{Row=19; Col=2; Old="13+1=14"; New="17+63=80"},
    @{Row=19; Col=3; Old="15+10=25"; New="65+21=86"},
    @{Row=19; Col=4; Old="15+67=82"; New="67-42=25"},
    @{Row=19; Col=5; Old="56-33=23"; New="30+1=31"},
    @{Row=20; Col=1; Old="86-69=17"; New="65+25=90"},
    @{Row=20; Col=2; Old="54+30=84"; New="16+71=87"},
    @{Row=20; Col=3; Old="75-34=41"; New="72-53=19"},
    @{Row=20; Col=4; Old="87-6=81"; New="79-1=78"},
    @{Row=20; Col=5; Old="20+75=95"; New="11+41=52"}
)

foreach ($item in $map) {
    $cell = $t.Cell($item.Row, $item.Col)
    $current = $cell.Range.Text.TrimEnd([char]13, [char]7)
    if ($current -ne $item.Old) {
        throw "Cell ($($item.Row),$($item.Col)) expected '$($item.Old)' but found '$current'"
    }
    $cell.Range.Text = $item.New
}

Write-Host "Applied $($map.Count) replacements"
